# Bug fix in Eduati data files:
# - Delete rows 45:87 on Sheet1 (they only contained leftover index values in column A)
# - Update the active sheet / selection state to reflect Sheet1 being active instead of Sheet3

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Remove the stray rows 45-87 from Sheet1 (data now stops at row 44)
$ws1.Rows("45:87").Delete()

# Reflect the new selection / scroll position on Sheet1
$ws1.Range("E64").Select()

# Sheet1 becomes the active/selected tab instead of Sheet3
$ws1.Activate()
$ws3.Select()
$ws1.Select()
